$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-6
# from 2023-09-15 (45184) to 2023-09-16 (45185)
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
